# Update "Moment for treatment group" values on the "data" sheet (column D)
# Correlation of past test scores with current portfolio and examen measures
# only for treatment group.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("D5").Value  = 2.490706599329265
$ws.Range("D6").Value  = 0.06892690360037453
$ws.Range("D7").Value  = -0.3633355102425982
$ws.Range("D8").Value  = 0.2290779302062357
$ws.Range("D9").Value  = 2.513774459734022
$ws.Range("D10").Value = 0.295579317546473
$ws.Range("D11").Value = 2.515486904511018
$ws.Range("D12").Value = 0.3455551069842551
$ws.Range("D13").Value = 0.3855308841340331
$ws.Range("D14").Value = 0.2158901897456601
$ws.Range("D15").Value = 0.2375479266643894
$ws.Range("D16").Value = 0.1209527201462497
$ws.Range("D17").Value = -0.08396974945207074
$ws.Range("D18").Value = -0.01276729934991094
$ws.Range("D19").Value = 0.6212943514050525
$ws.Range("D20").Value = 0.3830234347913151
$ws.Range("D21").Value = 0.1789474868051397
$ws.Range("D22").Value = 0.02360350061983014
